$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 85: new "Réalisation" entry about attempting to add comments ---
$ws.Range("A85").Value = 44705
$ws.Range("B85").Value = "Réalisation"
$ws.Range("C85").Value = 1
$ws.Range("D85").Value = "Tentative d'ajout de commentaires"
$ws.Range("E85").Value = "Ne fonctionne pas du à une erreur FK SQL"

# --- Row 86: documentation of the comment tests ---
$ws.Range("A86").Value = 44705
$ws.Range("B86").Value = "Réalisation"
$ws.Range("C86").Value = 0.25
$ws.Range("D86").Value = "Documentation des tests des commentaires"

# Copy the date formatting (column A, date number format + wrap text) from the
# last existing row down through the new rows, including the trailing blank
# row 87 that Excel leaves styled after the last entry.
$ws.Range("A84").Copy()
$ws.Range("A85:A87").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Grow the table ("Tableau1") so it covers the newly added rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F87"))

# Put the selection on the new first empty row, matching where Excel would
# leave the cursor after typing in the journal.
$ws.Range("A87").Select()
